$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H138").Value = 2409.8572
$ws.Range("I138").Value = 1130.2354
$ws.Range("J138").Value = 3280
$ws.Range("K138").Value = 3390.7062
$ws.Range("L138").Value = 9840
$ws.Range("M138").Value = 1749.2938
$ws.Range("N138").Value = -20120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10318.719
$ws.Range("I32").Value = 3264.262
$ws.Range("J32").Value = 23786.318
$ws.Range("K32").Value = 3264.262
$ws.Range("L32").Value = 23786.318
$ws.Range("M32").Value = -2977.262
$ws.Range("N32").Value = -24360.318

$ws.Range("H122").Value = 3403.5952
$ws.Range("I122").Value = 2864.4333
$ws.Range("J122").Value = 4751.5
$ws.Range("K122").Value = 8593.2999
$ws.Range("L122").Value = 14254.5
$ws.Range("M122").Value = -6143.2999
$ws.Range("N122").Value = -19154.5

$ws.Range("H141").Value = 45800
$ws.Range("J141").Value = 45800
$ws.Range("L141").Value = 45800
$ws.Range("N141").Value = -56160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 26000
$ws.Range("J50").Value = 26000
$ws.Range("L50").Value = 26000
$ws.Range("N50").Value = -27148

$ws.Range("H80").Value = 95.333336
$ws.Range("I80").Value = 92.25
$ws.Range("J80").Value = 96.454544
$ws.Range("K80").Value = 92.25
$ws.Range("L80").Value = 96.454544
$ws.Range("M80").Value = 905.75
$ws.Range("N80").Value = -2092.454544

$ws.Range("H83").Value = 95.333336
$ws.Range("I83").Value = 92.25
$ws.Range("J83").Value = 96.454544
$ws.Range("K83").Value = 461.25
$ws.Range("L83").Value = 482.27272
$ws.Range("M83").Value = 4530.75
$ws.Range("N83").Value = -10466.27272

$ws.Range("H99").Value = 1181.3572
$ws.Range("I99").Value = 1331.25
$ws.Range("J99").Value = 981.5
$ws.Range("K99").Value = 1331.25
$ws.Range("L99").Value = 981.5
$ws.Range("M99").Value = 166.75
$ws.Range("N99").Value = -3977.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5883569.5
$ws.Range("I31").Value = 9091620
$ws.Range("J31").Value = 4349284.5
$ws.Range("K31").Value = 9091620
$ws.Range("L31").Value = 4349284.5
$ws.Range("M31").Value = -9091325
$ws.Range("N31").Value = -4349874.5

$ws.Range("H34").Value = 5883569.5
$ws.Range("I34").Value = 9091620
$ws.Range("J34").Value = 4349284.5
$ws.Range("K34").Value = 9091620
$ws.Range("L34").Value = 4349284.5
$ws.Range("M34").Value = -9091418
$ws.Range("N34").Value = -4349688.5

$ws.Range("H48").Value = 5149.75
$ws.Range("J48").Value = 5149.75
$ws.Range("L48").Value = 5149.75
$ws.Range("N48").Value = -6101.75

$ws.Range("H68").Value = 10268
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 10268
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H99").Value = 5822.4
$ws.Range("I99").Value = 6032
$ws.Range("J99").Value = 5333.3335
$ws.Range("K99").Value = 6032
$ws.Range("L99").Value = 5333.3335
$ws.Range("M99").Value = -4534
$ws.Range("N99").Value = -8329.333500000001

$ws.Range("H122").Value = 1568.2354
$ws.Range("I122").Value = 728.9167
$ws.Range("J122").Value = 3582.6
$ws.Range("K122").Value = 2186.7501
$ws.Range("L122").Value = 10747.8
$ws.Range("M122").Value = 263.2498999999998
$ws.Range("N122").Value = -15647.8

$ws.Range("H126").Value = 5822.4
$ws.Range("I126").Value = 6032
$ws.Range("J126").Value = 5333.3335
$ws.Range("K126").Value = 18096
$ws.Range("L126").Value = 16000.0005
$ws.Range("M126").Value = -15626
$ws.Range("N126").Value = -20940.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 549.8125
$ws.Range("I34").Value = 124.625
$ws.Range("J34").Value = 975
$ws.Range("K34").Value = 373.875
$ws.Range("L34").Value = 2925
$ws.Range("M34").Value = -289.875
$ws.Range("N34").Value = -3093

$ws.Range("H39").Value = 2192.9167
$ws.Range("J39").Value = 2192.9167
$ws.Range("L39").Value = 6578.750100000001
$ws.Range("N39").Value = -7166.750100000001

$ws.Range("H55").Value = 1839.6364
$ws.Range("I55").Value = 387.2
$ws.Range("J55").Value = 3050
$ws.Range("K55").Value = 1161.6
$ws.Range("L55").Value = 9150
$ws.Range("M55").Value = -984.5999999999999
$ws.Range("N55").Value = -9504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3468.875
$ws.Range("I43").Value = 1291.8334
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 1291.8334
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = -1140.8334
$ws.Range("N43").Value = -10302

$ws.Range("H46").Value = 15000
$ws.Range("J46").Value = 15000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15312

$ws.Range("H57").Value = 10000
$ws.Range("I57").Value = 1000
$ws.Range("J57").Value = 19000
$ws.Range("K57").Value = 1000
$ws.Range("L57").Value = 19000
$ws.Range("M57").Value = -180
$ws.Range("N57").Value = -20640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2954.08
$ws.Range("I7").Value = 2875.389
$ws.Range("J7").Value = 3156.4285
$ws.Range("K7").Value = 2875.389
$ws.Range("L7").Value = 3156.4285
$ws.Range("M7").Value = -2763.389
$ws.Range("N7").Value = -3380.4285

$ws.Range("H40").Value = 3836.7297
$ws.Range("I40").Value = 3550.3103
$ws.Range("J40").Value = 4875
$ws.Range("K40").Value = 3550.3103
$ws.Range("L40").Value = 4875
$ws.Range("M40").Value = -3414.3103
$ws.Range("N40").Value = -5147

$ws.Range("H122").Value = 2494.647
$ws.Range("I122").Value = 2167
$ws.Range("J122").Value = 2724
$ws.Range("K122").Value = 6501
$ws.Range("L122").Value = 8172
$ws.Range("M122").Value = -4051
$ws.Range("N122").Value = -13072

$ws.Range("H126").Value = 2954.08
$ws.Range("I126").Value = 2875.389
$ws.Range("J126").Value = 3156.4285
$ws.Range("K126").Value = 8626.167000000001
$ws.Range("L126").Value = 9469.2855
$ws.Range("M126").Value = -6156.167000000001
$ws.Range("N126").Value = -14409.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 55556990
$ws.Range("I122").Value = 125001290
$ws.Range("J122").Value = 1555.9
$ws.Range("K122").Value = 375003870
$ws.Range("L122").Value = 4667.700000000001
$ws.Range("M122").Value = -375001420
$ws.Range("N122").Value = -9567.700000000001
